$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($Row, $Col, $Val) {
    $c = $ws.Cells.Item($Row, $Col)
    $c.NumberFormat = "@"
    $c.Value = $Val
}

# Row 2
$ws.Cells.Item(2, 4).Value = "30.732.33"
$ws.Cells.Item(2, 5).Value = "  +0.16%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.910.64"

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.11%  "

# Row 5
Set-Text 5 4 "239.30"
$ws.Cells.Item(5, 5).Value = "  -1.06%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.02%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.02%  "

# Row 8
Set-Text 8 4 "0.2957"
$ws.Cells.Item(8, 5).Value = "  +0.88%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +0.14%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "1.930.51"
$ws.Cells.Item(10, 5).Value = "  +1.98%  "

# Row 11
Set-Text 11 4 "17.06"
$ws.Cells.Item(11, 5).Value = "  -0.66%  "

# Row 12
Set-Text 12 4 "0.07361"
$ws.Cells.Item(12, 5).Value = "  +1.42%  "

# Row 13
Set-Text 13 4 "5.157"
$ws.Cells.Item(13, 5).Value = "  +2.78%  "

# Row 14
Set-Text 14 4 "88.27"
$ws.Cells.Item(14, 5).Value = "  -2.59%  "

# Row 15
Set-Text 15 4 "0.6697"
$ws.Cells.Item(15, 5).Value = "  -0.69%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "30.701.22"
$ws.Cells.Item(16, 5).Value = "  +0.08%  "

# Row 17
Set-Text 17 4 "0.000007891"
$ws.Cells.Item(17, 5).Value = "  -0.78%  "

# Row 18
Set-Text 18 4 "13.46"
$ws.Cells.Item(18, 5).Value = "  +2.93%  "

# Row 19
Set-Text 19 4 "1.002"
$ws.Cells.Item(19, 5).Value = "  +0.06%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "2.183.22"
$ws.Cells.Item(20, 5).Value = "  +1.95%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "Uniswap"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-Text 21 4 "5.305"
$ws.Cells.Item(21, 5).Value = "  +10.49%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "BinanceUSD"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-Text 22 4 "1.002"
$ws.Cells.Item(22, 5).Value = "  -0.01%  "

# Row 23
Set-Text 23 4 "195.16"
$ws.Cells.Item(23, 5).Value = "  +3.30%  "

# Row 24
Set-Text 24 4 "6.243"
$ws.Cells.Item(24, 5).Value = "  +2.70%  "

# Row 25
Set-Text 25 4 "9.634"
$ws.Cells.Item(25, 5).Value = "  +3.12%  "

# Row 26
Set-Text 26 4 "163.01"
$ws.Cells.Item(26, 5).Value = "  +3.85%  "

# Row 27
Set-Text 27 4 "18.58"
$ws.Cells.Item(27, 5).Value = "  -1.33%  "

# Row 28
Set-Text 28 4 "1.947"
$ws.Cells.Item(28, 5).Value = "  +3.02%  "

# Row 29
Set-Text 29 4 "1.475"
$ws.Cells.Item(29, 5).Value = "  +5.06%  "

# Row 30
Set-Text 30 4 "4.372"
$ws.Cells.Item(30, 5).Value = "  +2.47%  "

# Row 31
Set-Text 31 4 "0.09122"
$ws.Cells.Item(31, 5).Value = "  +0.55%  "

# Row 32
Set-Text 32 4 "4.052"
$ws.Cells.Item(32, 5).Value = "  +1.48%  "

# Row 33
Set-Text 33 4 "0.05250"
$ws.Cells.Item(33, 5).Value = "  +0.63%  "

# Row 34
Set-Text 34 4 "0.7379"
$ws.Cells.Item(34, 5).Value = "  +0.33%  "

# Row 35
Set-Text 35 4 "1.110"
$ws.Cells.Item(35, 5).Value = "  +0.55%  "

# Row 36
Set-Text 36 4 "2.725"
$ws.Cells.Item(36, 5).Value = "  -1.43%  "

# Row 37
Set-Text 37 4 "0.01818"
$ws.Cells.Item(37, 5).Value = "  -0.54%  "

# Row 38
Set-Text 38 4 "2.716"
$ws.Cells.Item(38, 5).Value = "  +1.35%  "

# Row 39
Set-Text 39 4 "0.9211"
$ws.Cells.Item(39, 5).Value = "  -0.65%  "

# Row 40
Set-Text 40 4 "2.069"
$ws.Cells.Item(40, 5).Value = "  -2.40%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +29.77%  "

# Row 42
Set-Text 42 4 "0.4435"
$ws.Cells.Item(42, 5).Value = "  +1.07%  "

# Row 43
Set-Text 43 4 "106.74"
$ws.Cells.Item(43, 5).Value = "  +1.66%  "

# Row 44
Set-Text 44 4 "5.902"
$ws.Cells.Item(44, 5).Value = "  +2.98%  "

# Row 45
Set-Text 45 4 "1.000"
$ws.Cells.Item(45, 5).Value = "  +0.06%  "

# Row 46
Set-Text 46 4 "0.1383"
$ws.Cells.Item(46, 5).Value = "  +2.70%  "

# Row 47
Set-Text 47 4 "7.582"
$ws.Cells.Item(47, 5).Value = "  +0.90%  "

# Row 48
Set-Text 48 4 "35.44"
$ws.Cells.Item(48, 5).Value = "  +5.48%  "

# Row 49
Set-Text 49 4 "9.059"
$ws.Cells.Item(49, 5).Value = "  +3.88%  "

# Row 50
Set-Text 50 4 "0.05858"
$ws.Cells.Item(50, 5).Value = "  -0.01%  "

# Row 51
Set-Text 51 4 "0.3992"
$ws.Cells.Item(51, 5).Value = "  +1.75%  "
